# Auto-generated script to apply Midgardsormr_Profits.xlsx market-data refresh
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ
# columns (H, I, J, K, L, M, N) on affected leve rows across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1973.2
$ws.Range("J17").Value = 2052.111
$ws.Range("L17").Value = 6156.333
$ws.Range("N17").Value = -6492.333
$ws.Range("H18").Value = 1199
$ws.Range("I18").Value = 1199
$ws.Range("K18").Value = 1199
$ws.Range("M18").Value = -915
$ws.Range("H112").Value = 6519.24
$ws.Range("I112").Value = 322.5
$ws.Range("K112").Value = 967.5
$ws.Range("M112").Value = 140.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1690.9183
$ws.Range("I2").Value = 1254.0278
$ws.Range("K2").Value = 1254.0278
$ws.Range("M2").Value = -1141.0278
$ws.Range("H32").Value = 3306.41
$ws.Range("I32").Value = 2996.299
$ws.Range("J32").Value = 13333.333
$ws.Range("K32").Value = 2996.299
$ws.Range("L32").Value = 13333.333
$ws.Range("M32").Value = -2709.299
$ws.Range("N32").Value = -13907.333
$ws.Range("H45").Value = 3921.9092
$ws.Range("J45").Value = 6068.1
$ws.Range("L45").Value = 6068.1
$ws.Range("N45").Value = -6822.1
$ws.Range("H61").Value = 1991.3334
$ws.Range("I61").Value = 1640.8868
$ws.Range("K61").Value = 1640.8868
$ws.Range("M61").Value = -1428.8868
$ws.Range("H74").Value = 190436.16
$ws.Range("I74").Value = 201195.4
$ws.Range("K74").Value = 201195.4
$ws.Range("M74").Value = -200321.4
$ws.Range("H77").Value = 190436.16
$ws.Range("I77").Value = 201195.4
$ws.Range("K77").Value = 1005977
$ws.Range("M77").Value = -1001609
$ws.Range("H102").Value = 4023.889
$ws.Range("I102").Value = 3030.3635
$ws.Range("K102").Value = 3030.3635
$ws.Range("M102").Value = -1408.3635
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H116").Value = 1690.9183
$ws.Range("I116").Value = 1254.0278
$ws.Range("K116").Value = 1254.0278
$ws.Range("M116").Value = 1039.9722
$ws.Range("H136").Value = 1991.3334
$ws.Range("I136").Value = 1640.8868
$ws.Range("K136").Value = 4922.6604
$ws.Range("M136").Value = -2372.6604
$ws.Range("H137").Value = 118450
$ws.Range("J137").Value = 118450
$ws.Range("L137").Value = 118450
$ws.Range("N137").Value = -128650

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1690.9183
$ws.Range("I3").Value = 1254.0278
$ws.Range("K3").Value = 1254.0278
$ws.Range("M3").Value = -1140.0278
$ws.Range("H94").Value = 1106.1428
$ws.Range("J94").Value = 998
$ws.Range("L94").Value = 998
$ws.Range("N94").Value = -1900
$ws.Range("H132").Value = 90149.664
$ws.Range("J132").Value = 90149.664
$ws.Range("L132").Value = 90149.664
$ws.Range("N132").Value = -100269.664
$ws.Range("H134").Value = 2177.9482
$ws.Range("I134").Value = 2026.9231
$ws.Range("K134").Value = 6080.7693
$ws.Range("M134").Value = -3545.7693

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("H62").Value = 14045.091
$ws.Range("I62").Value = 4999.2
$ws.Range("K62").Value = 4999.2
$ws.Range("M62").Value = -4375.2
$ws.Range("H65").Value = 14045.091
$ws.Range("I65").Value = 4999.2
$ws.Range("K65").Value = 24996
$ws.Range("M65").Value = -21876
$ws.Range("H94").Value = 1340.2222
$ws.Range("J94").Value = 1017.4286
$ws.Range("L94").Value = 1017.4286
$ws.Range("N94").Value = -1919.4286
$ws.Range("H134").Value = 2508.487
$ws.Range("I134").Value = 2044.5416
$ws.Range("J134").Value = 3250.8
$ws.Range("K134").Value = 6133.6248
$ws.Range("L134").Value = 9752.400000000001
$ws.Range("M134").Value = -3598.6248
$ws.Range("N134").Value = -14822.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3847.0667
$ws.Range("I34").Value = 661.75
$ws.Range("K34").Value = 1985.25
$ws.Range("M34").Value = -1901.25
$ws.Range("H39").Value = 4302.25
$ws.Range("J39").Value = 4753.357
$ws.Range("L39").Value = 14260.071
$ws.Range("N39").Value = -14848.071
$ws.Range("H103").Value = 731.3333
$ws.Range("J103").Value = 1449.5
$ws.Range("L103").Value = 4348.5
$ws.Range("N103").Value = -6106.5
$ws.Range("H114").Value = 1156.375
$ws.Range("J114").Value = 1498.3334
$ws.Range("L114").Value = 4495.0002
$ws.Range("N114").Value = -11003.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 30000
$ws.Range("J38").Value = 30000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30926
$ws.Range("H97").Value = 1743.2916
$ws.Range("I97").Value = 1477.5625
$ws.Range("J97").Value = 2274.75
$ws.Range("K97").Value = 1477.5625
$ws.Range("L97").Value = 2274.75
$ws.Range("M97").Value = -981.5625
$ws.Range("N97").Value = -3266.75
$ws.Range("H113").Value = 2094.6072
$ws.Range("I113").Value = 2116.9583
$ws.Range("J113").Value = 1960.5
$ws.Range("K113").Value = 2116.9583
$ws.Range("L113").Value = 1960.5
$ws.Range("M113").Value = 53.04170000000022
$ws.Range("N113").Value = -6300.5
$ws.Range("H132").Value = 2029.8
$ws.Range("I132").Value = 2099.75
$ws.Range("K132").Value = 6299.25
$ws.Range("M132").Value = -3769.25
$ws.Range("H140").Value = 185473.28
$ws.Range("J140").Value = 185473.28
$ws.Range("L140").Value = 185473.28
$ws.Range("N140").Value = -195833.28

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1504.2354
$ws.Range("I22").Value = 909.7143
$ws.Range("K22").Value = 909.7143
$ws.Range("M22").Value = -614.7143
$ws.Range("H27").Value = 1504.2354
$ws.Range("I27").Value = 909.7143
$ws.Range("K27").Value = 909.7143
$ws.Range("M27").Value = -802.7143
$ws.Range("H40").Value = 5498.533
$ws.Range("I40").Value = 5382.923
$ws.Range("K40").Value = 5382.923
$ws.Range("M40").Value = -5246.923
$ws.Range("H55").Value = 2925.4211
$ws.Range("I55").Value = 597.1429000000001
$ws.Range("K55").Value = 597.1429000000001
$ws.Range("M55").Value = -424.1429000000001
$ws.Range("H61").Value = 987.7368
$ws.Range("I61").Value = 939.9231
$ws.Range("K61").Value = 939.9231
$ws.Range("M61").Value = -737.9231
$ws.Range("H100").Value = 11908.909
$ws.Range("I100").Value = 15062.25
$ws.Range("K100").Value = 15062.25
$ws.Range("M100").Value = -14521.25
$ws.Range("H113").Value = 987.7368
$ws.Range("I113").Value = 939.9231
$ws.Range("K113").Value = 939.9231
$ws.Range("M113").Value = 1230.0769
$ws.Range("H132").Value = 3213.3447
$ws.Range("I132").Value = 3247.4546
$ws.Range("K132").Value = 9742.363799999999
$ws.Range("M132").Value = -7212.363799999999
$ws.Range("H137").Value = 54667
$ws.Range("J137").Value = 37000.5
$ws.Range("L137").Value = 37000.5
$ws.Range("N137").Value = -47200.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 326.21054
$ws.Range("I100").Value = 318.5
$ws.Range("K100").Value = 637
$ws.Range("M100").Value = -96
$ws.Range("H107").Value = 1722.25
$ws.Range("I107").Value = 2628
$ws.Range("K107").Value = 7884
$ws.Range("M107").Value = -5964
$ws.Range("H132").Value = 6118106.5
$ws.Range("I132").Value = 7377457.5
$ws.Range("J132").Value = 1257.8572
$ws.Range("K132").Value = 22132372.5
$ws.Range("L132").Value = 3773.5716
$ws.Range("M132").Value = -22129842.5
$ws.Range("N132").Value = -8833.571599999999
